$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert the two new rows.
#    - A new row is inserted at row 16 (becomes "GenerateScVsBulkDataQuantile.R")
#      pushing the old rows 16-24 down to 17-25.
#    - A second new row is inserted at row 19 (becomes "GenFigScVsBulkQuantile.R"),
#      i.e. right after the row that now holds "GenFigScVsBulk.R" (row 18),
#      pushing everything below further down.
# ---------------------------------------------------------------------------
$ws.Rows.Item(16).Insert()
$ws.Rows.Item(19).Insert()

# ---------------------------------------------------------------------------
# 2. Fill in the content of the two new rows.
#    Values are written in the same order the shared strings were added to
#    the workbook when it was originally authored (A16, A19, B19, then B16),
#    so new shared-string entries line up with the expected indices 39-42.
# ---------------------------------------------------------------------------
$ws.Range("A16").Value = 'GenerateScVsBulkDataQuantile.R'
$ws.Range("A19").Value = 'GenFigScVsBulkQuantile.R'
$ws.Range("B19").Value = 'This code is copied from GenFigScVsBulk.R and is very similar. A manual check that we loaded the right file was done, and the results look reasonable and a bit different from the TMM-normalized data. No other verification was deemed necessary.'
$ws.Range("B16").Value = 'TC004, found at the end of the code - the rest of the code is copied from GenerateScVsBulkData.R and is therefore tested as part of those tests. The data looks quantile normalized. No other verification was deemed necessary.'

# The inserted rows copied formatting from the row above (style 3 / indent), clear it
# so column A looks like the other top-level "file" rows (no extra style).
$ws.Range("A16").Style = "Normal"
$ws.Range("B16").Style = "Normal"
$ws.Range("B16").WrapText = $true

$ws.Range("A19").Style = "Normal"
$ws.Range("B19").Style = "Normal"
$ws.Range("B19").WrapText = $true

# ---------------------------------------------------------------------------
# 3. Row heights - set the explicit (custom) heights used in the final layout.
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 45
$ws.Rows.Item(3).RowHeight = 45
$ws.Rows.Item(4).RowHeight = 60
$ws.Rows.Item(7).RowHeight = 30
$ws.Rows.Item(10).RowHeight = 30
$ws.Rows.Item(11).RowHeight = 30
$ws.Rows.Item(13).RowHeight = 45
$ws.Rows.Item(16).RowHeight = 45
$ws.Rows.Item(19).RowHeight = 45
$ws.Rows.Item(20).RowHeight = 30
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 30

# Rows that go back to the (new) default row height - remove their explicit height.
$defaultHeightRows = @(5, 6, 8, 9, 12, 14, 15, 17, 18, 23, 24, 25, 26)
foreach ($r in $defaultHeightRows) {
    $ws.Rows.Item($r).AutoFit()
}

# ---------------------------------------------------------------------------
# 4. Update the selected cell shown when the workbook is reopened.
# ---------------------------------------------------------------------------
$ws.Range("B17").Select()
